$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1162.1052
$ws.Range("I43").Value = 1108.9166
$ws.Range("K43").Value = 1108.9166
$ws.Range("M43").Value = -1039.9166
$ws.Range("H98").Value = 1575.4615
$ws.Range("I98").Value = 1575.4615
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1575.4615
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -77.46149999999989
$ws.Range("H113").Value = 3955.5454
$ws.Range("I113").Value = 2044.8182
$ws.Range("K113").Value = 2044.8182
$ws.Range("M113").Value = 1209.1818
$ws.Range("H122").Value = 1575.4615
$ws.Range("I122").Value = 1575.4615
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4726.3845
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2276.3845
$ws.Range("H138").Value = 2977.373
$ws.Range("J138").Value = 3156.83
$ws.Range("L138").Value = 9470.49
$ws.Range("N138").Value = -19750.49

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 349.75
$ws.Range("I4").Value = 349.75
$ws.Range("K4").Value = 349.75
$ws.Range("M4").Value = -233.75
$ws.Range("H32").Value = 4176.189
$ws.Range("I32").Value = 2799.353
$ws.Range("J32").Value = 8431.862999999999
$ws.Range("K32").Value = 2799.353
$ws.Range("L32").Value = 8431.862999999999
$ws.Range("M32").Value = -2512.353
$ws.Range("N32").Value = -9005.862999999999
$ws.Range("H97").Value = 1011094.25
$ws.Range("I97").Value = 1406433.4
$ws.Range("J97").Value = 783
$ws.Range("K97").Value = 1406433.4
$ws.Range("L97").Value = 783
$ws.Range("M97").Value = -1405937.4
$ws.Range("N97").Value = -1775
$ws.Range("H124").Value = 13231.571
$ws.Range("J124").Value = 13231.571
$ws.Range("L124").Value = 13231.571
$ws.Range("N124").Value = -23051.571
$ws.Range("H132").Value = 2368.8462
$ws.Range("I132").Value = 1694.1765
$ws.Range("K132").Value = 5082.529500000001
$ws.Range("M132").Value = -2552.529500000001
$ws.Range("H133").Value = 20000
$ws.Range("J133").Value = 20000
$ws.Range("L133").Value = 20000
$ws.Range("N133").Value = -25060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7584624
$ws.Range("I94").Value = 12988369
$ws.Range("K94").Value = 12988369
$ws.Range("M94").Value = -12987918
$ws.Range("H105").Value = 3290723.2
$ws.Range("I105").Value = 3907515.2
$ws.Range("K105").Value = 3907515.2
$ws.Range("M105").Value = -3905768.2
$ws.Range("H110").Value = 84000
$ws.Range("J110").Value = 84000
$ws.Range("L110").Value = 84000
$ws.Range("N110").Value = -92180
$ws.Range("H134").Value = 3895.7097
$ws.Range("I134").Value = 1418.2222
$ws.Range("J134").Value = 7326.077
$ws.Range("K134").Value = 4254.6666
$ws.Range("L134").Value = 21978.231
$ws.Range("M134").Value = -1719.6666
$ws.Range("N134").Value = -27048.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 6029.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H58").Value = 1625.1177
$ws.Range("I58").Value = 1319.8846
$ws.Range("J58").Value = 2617.125
$ws.Range("K58").Value = 1319.8846
$ws.Range("L58").Value = 2617.125
$ws.Range("M58").Value = -1116.8846
$ws.Range("N58").Value = -3023.125
$ws.Range("H88").Value = 32247
$ws.Range("I88").Value = 10997
$ws.Range("K88").Value = 10997
$ws.Range("M88").Value = -10591
$ws.Range("H91").Value = 32247
$ws.Range("I91").Value = 10997
$ws.Range("K91").Value = 10997
$ws.Range("M91").Value = -9593
$ws.Range("H107").Value = 1593.85
$ws.Range("I107").Value = 1798.9412
$ws.Range("K107").Value = 1798.9412
$ws.Range("M107").Value = 121.0588
$ws.Range("H136").Value = 1625.1177
$ws.Range("I136").Value = 1319.8846
$ws.Range("J136").Value = 2617.125
$ws.Range("K136").Value = 3959.6538
$ws.Range("L136").Value = 7851.375
$ws.Range("M136").Value = -1409.6538
$ws.Range("N136").Value = -12951.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 16093.471
$ws.Range("I110").Value = 3713
$ws.Range("J110").Value = 17744.2
$ws.Range("K110").Value = 11139
$ws.Range("L110").Value = 53232.60000000001
$ws.Range("M110").Value = -7049
$ws.Range("N110").Value = -61412.60000000001
$ws.Range("H134").Value = 2005.8
$ws.Range("I134").Value = 1891.3077
$ws.Range("K134").Value = 5673.9231
$ws.Range("M134").Value = -603.9231
$ws.Range("H139").Value = 3596.4285
$ws.Range("I139").Value = 3916
$ws.Range("J139").Value = 2797.5
$ws.Range("K139").Value = 11748
$ws.Range("L139").Value = 8392.5
$ws.Range("M139").Value = -6608
$ws.Range("N139").Value = -18672.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 99638.2
$ws.Range("I22").Value = 141930.42
$ws.Range("K22").Value = 141930.42
$ws.Range("M22").Value = -141635.42
$ws.Range("H27").Value = 99638.2
$ws.Range("I27").Value = 141930.42
$ws.Range("K27").Value = 141930.42
$ws.Range("M27").Value = -141823.42
$ws.Range("H46").Value = 6218
$ws.Range("I46").Value = 4033.3333
$ws.Range("K46").Value = 4033.3333
$ws.Range("M46").Value = -3845.3333
$ws.Range("H61").Value = 4448046
$ws.Range("I61").Value = 6948373.5
$ws.Range("K61").Value = 6948373.5
$ws.Range("M61").Value = -6948171.5
$ws.Range("H93").Value = 18520254
$ws.Range("I93").Value = 20835056
$ws.Range("J93").Value = 1839
$ws.Range("K93").Value = 20835056
$ws.Range("L93").Value = 1839
$ws.Range("M93").Value = -20833808
$ws.Range("N93").Value = -4335
$ws.Range("H113").Value = 4448046
$ws.Range("I113").Value = 6948373.5
$ws.Range("K113").Value = 6948373.5
$ws.Range("M113").Value = -6946203.5
$ws.Range("H122").Value = 4911.9165
$ws.Range("I122").Value = 3218.6365
$ws.Range("J122").Value = 6344.6924
$ws.Range("K122").Value = 9655.9095
$ws.Range("L122").Value = 19034.0772
$ws.Range("M122").Value = -7205.9095
$ws.Range("N122").Value = -23934.0772
$ws.Range("H127").Value = 54666
$ws.Range("J127").Value = 54666
$ws.Range("L127").Value = 54666
$ws.Range("N127").Value = -64586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 83334340
$ws.Range("I81").Value = 83334340
$ws.Range("K81").Value = 166668680
$ws.Range("M81").Value = -166667619
$ws.Range("H84").Value = 83334340
$ws.Range("I84").Value = 83334340
$ws.Range("K84").Value = 833343400
$ws.Range("M84").Value = -833338096
